$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update the "Marking" row (row 11): number of marks per correct answer 3 -> 5
$ws.Range("B11").Value = 5

# Update the "Total" row (row 12): total correct marks 30 -> 50
$ws.Range("B12").Value = 50

# Update the correct/total marks display text (e.g. "30/84" -> "50/140")
$ws.Range("E12").Value = "50/140"
